# Mei-2016.xlsx — "10 years Finalization data"
#
# The daily-data table (header + 31 daily rows, A9:K40) that lives on the
# original "Data Harian - Table" sheet gets finalized onto a fresh
# worksheet ("Sheet1"), inserted right after the original sheet and left
# as the active tab. The decorative BMKG logo picture on the original
# sheet is removed.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new sheet right after the existing one so the tab order becomes
# "Data Harian - Table", "Sheet1" (Excel names a freshly added sheet
# "Sheet1" automatically).
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)

# Copy the finalized table (header row 9 + data rows 10-40) so it lands
# at A1:K32 on the new sheet.
$src = $ws1.Range("A9:K40")
$dst = $ws2.Range("A1")

$src.Copy()
$dst.PasteSpecial(-4163)   # xlPasteValues - values & shared strings
$src.Copy()
$dst.PasteSpecial(-4122)   # xlPasteFormats - border/alignment styles

# Leave the whole pasted table selected on the new sheet.
$ws2.Range("A1:K32").Select()

# Restore/leave the original sheet's selection over the table it donated,
# then drop its embedded logo picture.
$ws1.Activate()
$ws1.Range("A9:K40").Select()
if ($ws1.Shapes.Count -gt 0) {
    for ($i = $ws1.Shapes.Count; $i -ge 1; $i--) {
        $ws1.Shapes.Item($i).Delete()
    }
}

# The finalized sheet is the one left on top when the file is (re)opened.
$ws2.Activate()
